$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision has been made to select `"Barbie`" for the assembly on Friday.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("D3").Value = "no_decision, "
$ws.Range("C4").Value = "MSG: None`n`nMSG: The conversation ended without a clear decision about which movie would be shown on Friday. Therefore, the appropriate action is to acknowledge that no decision can be made.`n"
$ws.Range("D4").Value = "no_decision, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday ended without a clear choice.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been recorded with the selection of `"Barbie`" for acquisition.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: I have recorded the decision as no decision regarding the movie for Friday was made.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday resulted in no agreement.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: No movie decision has been made.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: Since the committee did not reach a consensus on what movie to show on Friday, I will call the function to indicate that no decision was made.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been recorded; `"Oppenheimer`" will be the movie shown on Friday.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision process concluded without arriving at a definitive choice for Friday's movie.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded successfully.`n"
$ws.Range("D16").Value = "both_movies, "
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision to show a movie on Friday has not been made.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision process has concluded without reaching an agreement on a movie for Friday, so no action will be taken regarding acquiring movie rights.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has concluded without an agreement.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday's showing.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no definitive movie choice was made for Friday.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: No decision was made about which movie to show on Friday.`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights to `"Barbie`" will be acquired for the movie to be shown on Friday.`n"
$ws.Range("D25").Value = "Barbie_was_selected, "
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday was not reached.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision process has concluded with no movie chosen for Friday.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie.`"`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was made regarding the movie to be shown on Friday.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been recorded with no selection made for the movie to be shown on Friday.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies has been recorded.`n"
$ws.Range("D31").Value = "both_movies, "
$ws.Range("C32").Value = "MSG: None`n`nMSG: The committee did not reach a decision regarding the movie to be shown on Friday.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been confirmed.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has not been made.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision being made regarding the movie to show on Friday.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected for Friday.`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected.`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: No decision was made about the movie to be shown on Friday.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the selection of a movie for Friday's showing.`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie has been selected for Friday.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was made regarding the movie for Friday.`n"
$ws.Range("D44").Value = "no_decision, "
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not reached.`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: I have recorded the decision as no decision on which movie to show on Friday.`n"
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie selection has been made.`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be acquired for Friday.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been recorded successfully.`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" has been selected as the movie for Friday's event.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision about the movie to play on Friday cannot be made at this time.`n"
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" will be acquired for showing on Friday.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: No decision was made about the movie to show on Friday.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("D55").Value = "both_movies, "
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision remains that no movie will be shown on Friday.`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision indicates that no movie was selected for Friday's showing.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision has been made that no movie will be selected for Friday.`n"
$ws.Range("C60").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision about Friday's movie was not made.`n"
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday could not be made, resulting in a no decision.`n"
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded successfully.`n"
$ws.Range("C65").Value = "MSG: None`n`nMSG: The committee did not reach a decision regarding the movie to be shown on Friday.`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
